$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.9226
$ws.Range("B2").Value = 0.9772
$ws.Range("C2").Value = 0.2214285714285714
$ws.Range("D2").Value = 0.001219512195121951
$ws.Range("E2").Value = 0.998780487804878
$ws.Range("F2").Value = 0.7785714285714286
$ws.Range("G2").Value = 0.01441677588466579
$ws.Range("H2").Value = 0.998780487804878
$ws.Range("I2").Value = 0.1113240418118467
$ws.Range("J2").Value = 0.921919770773639
$ws.Range("K2").Value = 0.07808022922636104
$ws.Range("A3").Value = 0.9360000000000001
$ws.Range("B3").Value = 0.9429
$ws.Range("C3").Value = 0.4696035242290749
$ws.Range("D3").Value = 0.004286520022560632
$ws.Range("E3").Value = 0.9957134799774394
$ws.Range("F3").Value = 0.5303964757709252
$ws.Range("G3").Value = 0.06312292358803986
$ws.Range("H3").Value = 0.9957134799774394
$ws.Range("I3").Value = 0.2369450221258178
$ws.Range("J3").Value = 0.9361544172234595
$ws.Range("K3").Value = 0.06384558277654051
$ws.Range("A4").Value = 0.9258
$ws.Range("B4").Value = 0.9462
$ws.Range("C4").Value = 0.4011627906976744
$ws.Range("D4").Value = 0.01382694023193577
$ws.Range("E4").Value = 0.9861730597680642
$ws.Range("F4").Value = 0.5988372093023255
$ws.Range("G4").Value = 0.2006472491909385
$ws.Range("H4").Value = 0.9861730597680642
$ws.Range("I4").Value = 0.2074948654648051
$ws.Range("J4").Value = 0.9346861128725428
$ws.Range("K4").Value = 0.06531388712745723
$ws.Range("A5").Value = 0.9468
$ws.Range("B5").Value = 0.9064
$ws.Range("C5").Value = 0.7
$ws.Range("D5").Value = 0.02547274749721913
$ws.Range("E5").Value = 0.9745272525027808
$ws.Range("F5").Value = 0.3
$ws.Range("G5").Value = 0.7557755775577558
$ws.Range("H5").Value = 0.9745272525027808
$ws.Range("I5").Value = 0.3627363737486096
$ws.Range("J5").Value = 0.9665710503089144
$ws.Range("K5").Value = 0.03342894969108556
$ws.Range("A6").Value = 0.9091
$ws.Range("B6").Value = 0.9772999999999999
$ws.Range("C6").Value = 0.1527494908350306
$ws.Range("D6").Value = 0.008538478598358839
$ws.Range("E6").Value = 0.9914615214016411
$ws.Range("F6").Value = 0.8472505091649695
$ws.Range("G6").Value = 0.09254807692307693
$ws.Range("H6").Value = 0.9914615214016411
$ws.Range("I6").Value = 0.08064398471669471
$ws.Range("J6").Value = 0.9148674920699887
$ws.Range("K6").Value = 0.08513250793001126
$ws.Range("A7").Value = 0.8744
$ws.Range("B7").Value = 0.8002
$ws.Range("C7").Value = 0.9159192825112108
$ws.Range("D7").Value = 0.1296662274923145
$ws.Range("E7").Value = 0.8703337725076855
$ws.Range("F7").Value = 0.0840807174887892
$ws.Range("G7").Value = 15.74666666666667
$ws.Range("H7").Value = 0.8703337725076855
$ws.Range("I7").Value = 0.5227927550017626
$ws.Range("J7").Value = 0.990627343164209
$ws.Range("K7").Value = 0.009372656835791027
$ws.Range("A8").Value = 0.9414
$ws.Range("B8").Value = 0.8673999999999999
$ws.Range("C8").Value = 0.8862212943632568
$ws.Range("D8").Value = 0.05275381552753815
$ws.Range("E8").Value = 0.9472461844724619
$ws.Range("F8").Value = 0.1137787056367432
$ws.Range("G8").Value = 4.376146788990826
$ws.Range("H8").Value = 0.9472461844724619
$ws.Range("I8").Value = 0.4694875549453975
$ws.Range("J8").Value = 0.9874337099377449
$ws.Range("K8").Value = 0.01256629006225507
$ws.Range("A9").Value = 0.8435
$ws.Range("B9").Value = 0.7519
$ws.Range("C9").Value = 0.9455252918287937
$ws.Range("D9").Value = 0.1681899242086491
$ws.Range("E9").Value = 0.8318100757913509
$ws.Range("F9").Value = 0.05447470817120625
$ws.Range("G9").Value = 26.94642857142857
$ws.Range("H9").Value = 0.8318100757913509
$ws.Range("I9").Value = 0.5568576080187214
$ws.Range("J9").Value = 0.9925522010905705
$ws.Range("K9").Value = 0.007447798909429482
$ws.Range("A10").Value = 0.9146
$ws.Range("B10").Value = 0.8902
$ws.Range("C10").Value = 0.6252566735112937
$ws.Range("D10").Value = 0.05417682251274097
$ws.Range("E10").Value = 0.9458231774872591
$ws.Range("F10").Value = 0.3747433264887063
$ws.Range("G10").Value = 1.33972602739726
$ws.Range("H10").Value = 0.9458231774872591
$ws.Range("I10").Value = 0.3397167480120173
$ws.Range("J10").Value = 0.9589979779824759
$ws.Range("K10").Value = 0.04100202201752412
$ws.Range("A11").Value = 0.8915999999999999
$ws.Range("B11").Value = 0.9403
$ws.Range("C11").Value = 0.2586719524281467
$ws.Range("D11").Value = 0.0373707040373707
$ws.Range("E11").Value = 0.9626292959626293
$ws.Range("F11").Value = 0.7413280475718533
$ws.Range("G11").Value = 0.4491978609625669
$ws.Range("H11").Value = 0.9626292959626293
$ws.Range("I11").Value = 0.1480213282327587
$ws.Range("J11").Value = 0.9204509199191747
$ws.Range("K11").Value = 0.07954908008082529
